$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.00036019
$ws.Range("F2").Value = 0.022946815
$ws.Range("G2").Value = 0.0005916887039552061

$ws.Range("E3").Value = 0.091220866
$ws.Range("F3").Value = 0.105217551
$ws.Range("G3").Value = 0.09548243956603773
